# Add 2022-Q3 data: insert a new quarter sheet right after the summary
# sheet ("总计"), and add the corresponding summary row there too. All
# of the other quarter sheets keep their own name+data pair and simply
# shift one tab to the right to make room.

$wb = $excel.ActiveWorkbook

function Set-TextValue($rng, [string]$val) {
    # Force the cell to stay text-typed even though the value looks
    # numeric (matches the source data, which stores these as text).
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

# ---------------------------------------------------------------------
# 1) "总计" (summary) sheet: insert a new row for 2022-Q3 right after the
#    header row, pushing the existing quarters down by one row.
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item(1)

$summary.Rows.Item(2).Insert()
$summary.Range("B2:D2").ClearFormats()

# Copy the style of the index column from the row below (the previously
# existing data rows) so the new row matches the others.
$summary.Range("A3").Copy()
$summary.Range("A2").PasteSpecial(-4122)   # xlPasteFormats

$summary.Range("A2").Value = 0
$summary.Range("B2").Value = "2022-Q3"
$summary.Range("C2").Value = 2
$summary.Range("D2").Value = 2.25

# Renumber the index column (0-based) for the rows that shifted down.
$summary.Range("A3").Value = 1
$summary.Range("A4").Value = 2
$summary.Range("A5").Value = 3
$summary.Range("A6").Value = 4
$summary.Range("A7").Value = 5

# ---------------------------------------------------------------------
# 2) Insert a brand-new "2022-Q3" fund-detail sheet right after "总计".
#    Clone the existing "2022-Q2" sheet (position 2) so the new sheet
#    keeps identical layout/styling, place it before that sheet (so it
#    lands in slot 2 and "2022-Q2" simply slides down to slot 3), then
#    overwrite it with the new quarter's figures.
# ---------------------------------------------------------------------
$templateSheet = $wb.Worksheets.Item(2)
$templateSheet.Copy($templateSheet, $null)

$q3 = $wb.Worksheets.Item(2)
$q3.Name = "2022-Q3"

# Row 2 - 008763 天弘越南市场股票（QDII）A
Set-TextValue $q3.Range("D2") "20.44"
Set-TextValue $q3.Range("E2") "90.19"
Set-TextValue $q3.Range("F2") "6.34"
Set-TextValue $q3.Range("G2") "1.2959"
$q3.Range("H2").Value = 4

# Row 3 - 008764 天弘越南市场股票（QDII）C
Set-TextValue $q3.Range("D3") "15.02"
Set-TextValue $q3.Range("E3") "90.19"
Set-TextValue $q3.Range("F3") "6.34"
Set-TextValue $q3.Range("G3") "0.9523"
$q3.Range("H3").Value = 4

# ---------------------------------------------------------------------
# 3) Restore the "selected tab" marker onto the last sheet (2021-Q2),
#    which is where it lived before the new sheet was inserted.
# ---------------------------------------------------------------------
$lastIndex = $wb.Worksheets.Count
$wb.Worksheets.Item($lastIndex).Select()
